# ------------------------------------------------------------------
# 221122 suhwa new crawling_ing
# Adds a new "H" notes column, fills in new status-tracking content
# for the 11/22 ~ 11/23 work columns (E/F), widens the updated
# columns, and adjusts a couple of row heights.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Move the long "notes" text out of G2 into the new H column ---
$notesText = $ws.Range("G2").Value()
$ws.Range("G2:G9").UnMerge()
$ws.Range("G2").ClearContents()
$ws.Range("G6:G8").Clear()
$ws.Range("H2").Value = $notesText
$ws.Range("H2:H9").Merge()

# --- 2. Fill in the new progress-report cells ---
$ws.Range("E2").Value = "댓글 마무리`nui 바꾸기"
$ws.Range("F2").Value = "유지보수`nui 바꾸기`n그래프"

$ws.Range("E4").Value = "관심지역, 상권정보"
$ws.Range("F4").Value = "상권정보 마무리`n뉴스"

$ws.Range("E5").Value = "관심지역`n병원`n코로나 진료소`n지도 뿌리기"
$ws.Range("F5").Value = "상권정보`n유지보수`n뉴스 시작"

# --- 3. New H1 date cell (11/24 -> serial 44889) ---
$ws.Range("H1").Value = 44889

# --- 4. Re-apply formatting using existing, already-styled cells as
#        templates, so the engine reuses/derives matching cellXfs. ---

# 4a. Date format for H1 (same look as C1:G1)
$ws.Range("C1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# 4b. G2 / G3 drop the old "notes" wrap style and become plain
#     centered (same look as C4/D3)
$ws.Range("C4").Copy()
$ws.Range("G2:G3").PasteSpecial(-4122)

# 4c. E2/F2/E5/F5/F4/G4/G5/G9 take on the centered wrap-text look
#     (same look as C2/D2)
$ws.Range("C2").Copy()
$ws.Range("E2:F2").PasteSpecial(-4122)
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("E5:F5").PasteSpecial(-4122)
$ws.Range("G4").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G9").PasteSpecial(-4122)

# 4d. E4 keeps the plain centered (no-wrap) look already used by C4/D4
$ws.Range("C4").Copy()
$ws.Range("E4").PasteSpecial(-4122)

# 4e. New H column (notes) takes the centered wrap-text look too
$ws.Range("C2").Copy()
$ws.Range("H2:H9").PasteSpecial(-4122)

# 4f. B9 switches from the C2-style wrap font to the B8-style wrap font
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 5. Column widths for the edited/added columns ---
$ws.Columns("E").ColumnWidth = 13.86
$ws.Columns("F").ColumnWidth = 11.71
$ws.Columns("G").ColumnWidth = 11.71
$ws.Columns("H").ColumnWidth = 39.71

# --- 6. Row height tweaks ---
$ws.Rows(2).RowHeight = 37.5
$ws.Rows(4).RowHeight = 24
$ws.Rows(5).RowHeight = 48

# --- 7. Restore the active selection shown in the workbook ---
$ws.Range("E3").Select()
